$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit swaps the data of row 3 and row 4 (two species observations
# traded places in the export), except for the few columns that already
# held identical values in both rows (C, D, I, T, U, V, W, Y, AA, AD, AE,
# AG, AT, AY) which are left untouched.

# --- Row 3 gets row 4's former record data ---
$ws.Range("A3").Value = 130861152
$ws.Range("B3").Value = 91804
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("P3").Value = "Djupbäcken, Jmt"
$ws.Range("Q3").Value = 442868
$ws.Range("R3").Value = 7039767
$ws.Range("S3").Value = 10

# Starttid/Sluttid (Z3/AB3) no longer apply to row 3 - clear them
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Publik kommentar (AC3) now carries row 4's former comment text.
# Set via a Text-formatted write then ClearFormats so it stays a plain
# literal string (no style/number-format side effects).
$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = "I stående levande gran med full längd."
$ws.Range("AC3").ClearFormats()

$ws.Range("AW3").Value = "Kristian Zackrisson"
$ws.Range("AX3").Value = "Kristian Zackrisson"

# --- Row 4 gets row 3's former record data ---
$ws.Range("A4").Value = 130853761
$ws.Range("B4").Value = 79244
$ws.Range("E4").Value = 230405
$ws.Range("F4").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G4").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("P4").Value = "Djupbäcken, Djupbäcken, Jmt"
$ws.Range("Q4").Value = 442771
$ws.Range("R4").Value = 7039709
$ws.Range("S4").Value = 20

# Starttid/Sluttid (Z4/AB4) now carry row 3's former "11:05" times.
# Written Text-formatted then ClearFormats so they remain plain literal
# strings rather than being auto-parsed into time serial numbers.
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "11:05"
$ws.Range("Z4").ClearFormats()

$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "11:05"
$ws.Range("AB4").ClearFormats()

# Publik kommentar (AC4) no longer applies to row 4 - clear it
$ws.Range("AC4").ClearContents()

$ws.Range("AW4").Value = "Maria Danvind"
$ws.Range("AX4").Value = "Maria Danvind"
